$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the repeated Print_Area defined-name pattern by three more entries ---
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "='Weekly Time Record'!`$A`$1:`$K`$27")
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "='Weekly Time Record'!`$A`$1:`$K`$27")
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "='Weekly Time Record'!`$A`$1:`$K`$27")

# --- Fill in the previously-blank Wednesday (row 15) and Friday (row 17) time punches ---
$ws.Range("C15").Value = 0.541666666666667
$ws.Range("D15").Value = 0.833333333333333
$ws.Range("C17").Value = 0.583333333333333
$ws.Range("D17").Value = 0.666666666666667

# --- Remove the stray comment note that lived in L18 (also drops it from sharedStrings) ---
$ws.Range("L18").Clear()

# --- Move the active selection to D18 ---
$ws.Range("D18").Select() | Out-Null

# --- Slightly narrow the data-entry columns (B, C:G, H, I:K) ---
$ws.Columns.Item(2).ColumnWidth = 8.25
$ws.Range("C1:G1").EntireColumn.ColumnWidth = 6.5
$ws.Columns.Item(8).ColumnWidth = 7.25
$ws.Range("I1:K1").EntireColumn.ColumnWidth = 6.5
